{"js": "// Add the IntelliJ \"Main Build Version\" (IU 183.5912.21) right after\n// \"IDEA 2018.3.5\" in the \"Software environment:\" line, turning:\n//   \"Software environment: IntelliJ IDEA 2018.3.5, JUnit 5.4.2\"\n// into:\n//   \"Software environment: IntelliJ IDEA 2018.3.5 IU 183.5912.21, JUnit 5.4.2\"\n\nconst body = context.document.body;\n\nconst results = body.search(\"IDEA 2018.3.5\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"IDEA 2018.3.5\" in the document body.');\n}\n\n// Insert the new text immediately after the found range so it keeps the\n// same run formatting (bold, navy color, size 20) as the surrounding text.\nresults.items[0].insertText(\" IU 183.5912.21\", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# Add the IntelliJ \"Main Build Version\" (IU 183.5912.21) right after\n# \"IDEA 2018.3.5\" in the \"Software environment:\" line, turning:\n#   \"Software environment: IntelliJ IDEA 2018.3.5, JUnit 5.4.2\"\n# into:\n#   \"Software environment: IntelliJ IDEA 2018.3.5 IU 183.5912.21, JUnit 5.4.2\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"IDEA 2018.3.5\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found) {\n    $rng = $find.Parent\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\" IU 183.5912.21\")\n}\n"}
